# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.606.19"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.290.04"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D5").Value = "'304.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'95.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").Value = "'34.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").Value = "'0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "'18.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "2.644.62"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "2.288.26"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "'0.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "42.497.50"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "'5.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").Value = "'66.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'235.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'25.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'165.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'8.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'32.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'4.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'17.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.0685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.101"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.992.52"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0277"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'17.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").Value = "'9.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.09%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.52%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'53.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.510.80"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'70.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.46%  "
